$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "last row" (80) loses its special last-row date style and
# reverts to the regular date style used by all other data rows.
$ws.Range("A80").NumberFormat = $ws.Range("A79").NumberFormat

# Append the new daily update row (81) using the style that previously
# marked the last row (now moved down to row 81).
$ws.Range("A81").NumberFormat = "YYYY-MM-DD"
$ws.Range("A81").Value = 45821
$ws.Range("B81").Value = 344
$ws.Range("C81").Value = 346
$ws.Range("D81").Value = 350
